$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the old "_GoBack" bookmark that sat after the "Design"
#    paragraph (it is relocated to the end of the paragraph edited
#    below).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Append the "older Excel versions" sentence to the paragraph that
#    talks about selecting the 'ExchangeRate-AddIn64' file. The whole
#    paragraph is rebuilt (original runs + new runs) and pushed back
#    in with InsertXML so the newly typed word gets a spell-check
#    proofErr wrapper, exactly like Word would add for a fresh typo.
# ------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*ExchangeRate-AddIn64*") {
        $r = $p.Range

        $xml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' + `
            '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
            '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
            '<pkg:xmlData>' + `
            '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
            '<w:body>' + `
            '<w:p w:rsidR="00172B1C" w:rsidRDefault="00D555A0" w:rsidP="00172B1C">' + `
            '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr>' + `
            '<w:r><w:t xml:space="preserve">Click on </w:t></w:r>' + `
            '<w:r w:rsidR="00172B1C"><w:t>‘</w:t></w:r>' + `
            '<w:r><w:t>Go</w:t></w:r>' + `
            '<w:r w:rsidR="00172B1C"><w:t>’. ‘Add-in’ screen will be opened. Browse the location of add-in and select ‘</w:t></w:r>' + `
            '<w:r w:rsidR="00172B1C" w:rsidRPr="00172B1C"><w:t>ExchangeRate-AddIn64</w:t></w:r>' + `
            '<w:r w:rsidR="00172B1C"><w:t xml:space="preserve">’ file. </w:t></w:r>' + `
            '<w:r><w:t>In case of older versions of Excel, select ‘</w:t></w:r>' + `
            '<w:proofErr w:type="spellStart"/>' + `
            '<w:r><w:t>ExchangeRate-AddIn</w:t></w:r>' + `
            '<w:proofErr w:type="spellEnd"/>' + `
            '<w:r><w:t>’ file.</w:t></w:r>' + `
            '</w:p>' + `
            '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

        $r.InsertXML($xml)
        break
    }
}

# ------------------------------------------------------------------
# 3) Drop the "_GoBack" bookmark back in, now at the very end of the
#    paragraph just edited (after "... file." and before the
#    paragraph mark).
#
#    Quirk: adding a *collapsed* Bookmarks range exactly at a
#    paragraph's end position (i.e. right on the paragraph mark)
#    resets the bookmark to the top of the document, so a throwaway
#    marker character is used to hold the spot, the bookmark is
#    anchored next to it, and the marker is deleted afterwards.
# ------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*In case of older versions*") {
        $endPos = $p.Range.End - 1

        $marker = $d.Range($endPos, $endPos)
        $marker.InsertAfter("X")

        $bmRange = $d.Range($endPos, $endPos)
        $d.Bookmarks.Add("_GoBack", $bmRange)

        $markerRange = $d.Range($endPos, $endPos + 1)
        $markerRange.Delete()
        break
    }
}
